$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# 1. Add LinkedIn to contact line
Replace-Text `
    "Pittsburgh, PA 15216 | (517) 515-1699 | quincy120@gmail.com | slimeq.github.io | github.com/SlimeQ" `
    "Pittsburgh, PA 15216 | (517) 515-1699 | quincy120@gmail.com | slimeq.github.io | github.com/SlimeQ | linkedin.com/in/quincy-campbell-131559b2"

# 2. Tin Drum bullet: warehouse-scale AR apps
Replace-Text `
    "Built warehouse-scale Unity AR apps for Magic Leap 2 (Android) for live performances (including KAGAMI)." `
    "Built Unity applications for warehouse-scale AR experiences on Magic Leap 2."

# 3. Tin Drum bullet: Android fleet tooling
Replace-Text `
    "Developed supporting platform: Android fleet tools (200+ devices), C#/ASP.NET Core/Blazor systems for metrics, file distribution, and AV sync." `
    "Built Android fleet tooling for 200+ devices."

# 4. Tin Drum bullet: C#/ASP.NET Core/Blazor services
Replace-Text `
    "Created reusable Unity plugins; implemented integration patterns (HTTP/WebSockets) and networking (Mirror/NGO)." `
    "Built C#/ASP.NET Core/Blazor services for metrics and content distribution."

# 5. Tin Drum bullet: volumetric capture pipeline wording
Replace-Text `
    "Implemented volumetric capture compression/decompression + playback pipeline." `
    "Implemented volumetric capture compression/decompression and playback pipeline."

# 6. Chameleon Power bullet: product visualization apps
Replace-Text `
    "Shipped Unity/PlayCanvas visualization products across PC, WebGL, iOS, Android." `
    "Shipped product visualization apps across PC, WebGL, iOS, and Android (Unity/PlayCanvas)."

# 7. Chameleon Power bullet: reusable Unity plugins wording
Replace-Text `
    "Created reusable Unity plugins/internal tools; optimized performance and build size for low-end hardware and WebGL." `
    "Built reusable Unity plugins/internal tools; optimized performance and build size for low-end hardware and WebGL."
